# Updates cryptos list prices/volumes (scraped data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.066.64"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "2.303.96"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").Value = "  +0.05%  "
$origStyle_D5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.85"
$ws.Range("D5").Style = $origStyle_D5
$ws.Range("E5").Value = "  -0.08%  "
$origStyle_D6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.91"
$ws.Range("D6").Style = $origStyle_D6
$ws.Range("E6").Value = "  -1.20%  "
$ws.Range("E7").Value = "  +4.12%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +1.00%  "
$origStyle_D10 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.74"
$ws.Range("D10").Style = $origStyle_D10
$ws.Range("E10").Value = "  -0.57%  "
$origStyle_D11 = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0789"
$ws.Range("D11").Style = $origStyle_D11
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("E12").Value = "  -0.03%  "
$origStyle_D13 = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.99"
$ws.Range("D13").Style = $origStyle_D13
$ws.Range("E13").Value = "  +0.74%  "
$origStyle_D14 = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.88"
$ws.Range("D14").Style = $origStyle_D14
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("D15").Value = "2.662.83"
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("D16").Value = "2.304.05"
$ws.Range("E16").Value = "  -1.11%  "
$origStyle_D17 = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.785"
$ws.Range("D17").Style = $origStyle_D17
$ws.Range("E17").Value = "  -1.64%  "
$ws.Range("D18").Value = "42.970.24"
$ws.Range("E18").Value = "  +0.68%  "
$origStyle_D19 = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.34"
$ws.Range("D19").Style = $origStyle_D19
$ws.Range("E19").Value = "  +7.85%  "
$ws.Range("D20").Value = "0.0₃0907"
$ws.Range("E20").Value = "  +0.87%  "
$origStyle_D21 = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.12"
$ws.Range("D21").Style = $origStyle_D21
$ws.Range("E21").Value = "  -1.01%  "
$origStyle_D22 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.28"
$ws.Range("D22").Style = $origStyle_D22
$ws.Range("E22").Value = "  +0.67%  "
$origStyle_D23 = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.53"
$ws.Range("D23").Style = $origStyle_D23
$ws.Range("E23").Value = "  +1.18%  "
$ws.Range("E24").Value = "  -1.42%  "
$origStyle_D26 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.42"
$ws.Range("D26").Style = $origStyle_D26
$ws.Range("E26").Value = "  -0.83%  "
$origStyle_D27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.70"
$ws.Range("D27").Style = $origStyle_D27
$ws.Range("E27").Value = "  +0.53%  "
$origStyle_D28 = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "167.88"
$ws.Range("D28").Style = $origStyle_D28
$ws.Range("E28").Value = "  -0.32%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  -0.99%  "
$origStyle_D31 = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.94"
$ws.Range("D31").Style = $origStyle_D31
$ws.Range("E31").Value = "  -4.46%  "
$origStyle_D32 = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.17"
$ws.Range("D32").Style = $origStyle_D32
$ws.Range("E32").Value = "  +4.10%  "
$origStyle_D33 = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").Style = $origStyle_D33
$origStyle_D34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.79"
$ws.Range("D34").Style = $origStyle_D34
$ws.Range("E34").Value = "  +4.57%  "
$origStyle_D35 = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.08"
$ws.Range("D35").Style = $origStyle_D35
$ws.Range("E35").Value = "  +3.37%  "
$ws.Range("E36").Value = "  +0.06%  "
$origStyle_D37 = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0689"
$ws.Range("D37").Style = $origStyle_D37
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("E38").Value = "  +1.28%  "
$ws.Range("E39").Value = "  +0.98%  "
$ws.Range("E40").Value = "  +2.26%  "
$ws.Range("E41").Value = "  -2.49%  "
$ws.Range("D42").Value = "2.011.34"
$ws.Range("E42").Value = "  +0.94%  "
$origStyle_D43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0289"
$ws.Range("D43").Style = $origStyle_D43
$ws.Range("E43").Value = "  +0.60%  "
$origStyle_D44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.18"
$ws.Range("D44").Style = $origStyle_D44
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("E45").Value = "  -3.89%  "
$origStyle_D46 = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.41"
$ws.Range("D46").Style = $origStyle_D46
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("E47").Value = "  -1.77%  "
$origStyle_D48 = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.42"
$ws.Range("D48").Style = $origStyle_D48
$ws.Range("E48").Value = "  -2.23%  "
$ws.Range("D49").Value = "2.528.89"
$origStyle_D50 = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.53"
$ws.Range("D50").Style = $origStyle_D50
$ws.Range("E50").Value = "  +0.66%  "
$origStyle_D51 = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.34"
$ws.Range("D51").Style = $origStyle_D51
$ws.Range("E51").Value = "  +4.73%  "

Write-Output "Updated cryptos list"
